{"js": "// Remove the block of paragraphs that runs from right after the\n// \"Marine coastal biogeochemistry...\" bio paragraph through the end of\n// the \"Keywords: ...\" paragraph (inclusive). This deletes the blank\n// spacer paragraphs, the \"Author Contribution Statement:\" heading, the\n// long author-contribution paragraph, another spacer, and the Keywords\n// paragraph -- while leaving the final blank paragraph (and the page\n// break paragraph after it) untouched, per the commit \"move some\n// submission info to ms\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet startIndex = -1; // first paragraph to delete (right after the bio paragraph)\nlet endIndex = -1;   // last paragraph to delete (the \"Keywords\" paragraph)\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (startIndex === -1 && text.indexOf(\"Marine coastal biogeochemistry\") !== -1) {\n    startIndex = i + 1;\n  }\n  if (text.indexOf(\"Keywords\") === 0) {\n    endIndex = i;\n  }\n}\n\nif (startIndex !== -1 && endIndex !== -1 && endIndex >= startIndex) {\n  // Delete from the end backwards so earlier indices stay valid.\n  for (let i = endIndex; i >= startIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the block of paragraphs that runs from right after the\n# \"Marine coastal biogeochemistry...\" bio paragraph through the end of\n# the \"Keywords: ...\" paragraph (inclusive). This deletes the blank\n# spacer paragraphs, the \"Author Contribution Statement:\" heading, the\n# long author-contribution paragraph, another spacer, and the Keywords\n# paragraph -- while leaving the final blank paragraph (and the page\n# break paragraph after it) untouched, per the commit \"move some\n# submission info to ms\".\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$startIndex = -1\n$endIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text\n    if ($startIndex -eq -1 -and $t -like \"*Marine coastal biogeochemistry*\") {\n        $startIndex = $i + 1\n    }\n    if ($t -like \"Keywords*\") {\n        $endIndex = $i\n    }\n}\n\nif ($startIndex -ne -1 -and $endIndex -ne -1 -and $endIndex -ge $startIndex) {\n    $startPara = $paras.Item($startIndex)\n    $endPara = $paras.Item($endIndex)\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
